$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("R1").Range("G2").Value = "3883:24:06"
$wb.Worksheets.Item("R1").Range("G3").Value = "22:56:44"

$wb.Worksheets.Item("R2").Range("G2").Value = "12064:47:47"
$wb.Worksheets.Item("R2").Range("G3").Value = "3194:31:16"
$wb.Worksheets.Item("R2").Range("G4").Value = "432:42:50"

$wb.Worksheets.Item("R4").Range("G2").Value = "2910:37:36"
$wb.Worksheets.Item("R4").Range("G3").Value = "137:49:51"

$wb.Worksheets.Item("R5").Range("G2").Value = "384:36:35"

$wb.Worksheets.Item("R6").Range("G2").Value = "25:08:53"
